$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6666666666666666
$ws.Range("C2").Value = 0.4
$ws.Range("D2").Value = 0.5
$ws.Range("E2").Value = 5

$ws.Range("B3").Value = 0.6666666666666666
$ws.Range("C3").Value = 0.4
$ws.Range("D3").Value = 0.5
$ws.Range("E3").Value = 5

$ws.Range("B4").Value = 0.75
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 0.8571428571428571
$ws.Range("E4").Value = 6

$ws.Range("B5").Value = 0.6666666666666666
$ws.Range("C5").Value = 0.6666666666666666
$ws.Range("D5").Value = 0.6666666666666666
$ws.Range("E5").Value = 3

$ws.Range("B6").Value = 0.65
$ws.Range("C6").Value = 0.9285714285714286
$ws.Range("D6").Value = 0.7647058823529412
$ws.Range("E6").Value = 14

$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 2

$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 0.6666666666666666
$ws.Range("D8").Value = 0.8
$ws.Range("E8").Value = 3

$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 0.8
$ws.Range("D9").Value = 0.888888888888889
$ws.Range("E9").Value = 5

$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 0.6666666666666666
$ws.Range("D10").Value = 0.8
$ws.Range("E10").Value = 3

$ws.Range("B11").Value = 0.5
$ws.Range("C11").Value = 0.5
$ws.Range("D11").Value = 0.5
$ws.Range("E11").Value = 2

$ws.Range("B12").Value = 0.625
$ws.Range("C12").Value = 0.8333333333333334
$ws.Range("D12").Value = 0.7142857142857143
$ws.Range("E12").Value = 6

$ws.Range("B13").Value = 0.6666666666666666
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 0.8
$ws.Range("E13").Value = 2

$ws.Range("B14").Value = 1
$ws.Range("C14").Value = 0.6666666666666666
$ws.Range("D14").Value = 0.8
$ws.Range("E14").Value = 3

$ws.Range("B15").Value = 1
$ws.Range("C15").Value = 0.6666666666666666
$ws.Range("D15").Value = 0.8
$ws.Range("E15").Value = 3

$ws.Range("B16").Value = 0.8
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 0.888888888888889
$ws.Range("E16").Value = 4

$ws.Range("B17").Value = 1
$ws.Range("C17").Value = 0.5
$ws.Range("D17").Value = 0.6666666666666666
$ws.Range("E17").Value = 2

$ws.Range("B18").Value = 0
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 2

$ws.Range("B19").Value = 1
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 3

$ws.Range("B20").Value = 0
$ws.Range("C20").Value = 0
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 2

$ws.Range("B21").Value = 0
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 2

$ws.Range("B22").Value = 0.6
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 0.7499999999999999
$ws.Range("E22").Value = 3

$ws.Range("B23").Value = 0.6666666666666666
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 0.8
$ws.Range("E23").Value = 4

$ws.Range("B24").Value = 1
$ws.Range("C24").Value = 0.6666666666666666
$ws.Range("D24").Value = 0.8
$ws.Range("E24").Value = 3

$ws.Range("B25").Value = 0.735632183908046
$ws.Range("C25").Value = 0.735632183908046
$ws.Range("D25").Value = 0.735632183908046
$ws.Range("E25").Value = 0.735632183908046

$ws.Range("B26").Value = 0.7068840579710145
$ws.Range("C26").Value = 0.6679089026915114
$ws.Range("D26").Value = 0.6650976332562011
$ws.Range("E26").Value = 87

$ws.Range("B27").Value = 0.7243295019157088
$ws.Range("C27").Value = 0.735632183908046
$ws.Range("D27").Value = 0.7071010013200683
$ws.Range("E27").Value = 87
